# Updates the cryptocurrency price/volume table on Sheet1 to reflect the
# latest scrape (GitHub Actions refresh job). Values in columns D (Price)
# and E (Volume(1h)) are stored as text in the workbook (not numbers), and
# several look like numbers (e.g. "1.00", "629.61") or contain two dots as
# thousands separators (e.g. "98.714.72"). Plain COM ".Value" assignment
# would make Excel auto-convert numeric-looking strings into real numbers
# and strip formatting (e.g. "1.00" -> 1), so each target cell is forced to
# Text format before the value is written, then its style is reset back to
# "Normal" so no stray per-cell formatting is left behind.
# Rows 37/38 also swap their Coin/Link/Price/Volume data (RenderToken and
# Kaspa traded places in the ranking).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue "D2" "98.714.72"
Set-TextValue "E2" "  +0.32%  "
Set-TextValue "D3" "3.313.93"
Set-TextValue "E3" "  -1.77%  "
Set-TextValue "E4" "  +0.01%  "
Set-TextValue "D5" "255.43"
Set-TextValue "E5" "  +0.18%  "
Set-TextValue "D6" "629.61"
Set-TextValue "E6" "  +1.07%  "
Set-TextValue "D7" "1.47"
Set-TextValue "E7" "  +21.68%  "
Set-TextValue "D8" "0.413"
Set-TextValue "E8" "  +7.20%  "
Set-TextValue "E9" "  -0.01%  "
Set-TextValue "E10" "  +24.58%  "
Set-TextValue "D11" "3.311.41"
Set-TextValue "E11" "  -1.71%  "
Set-TextValue "E12" "  +2.92%  "
Set-TextValue "D13" "41.72"
Set-TextValue "E13" "  +16.33%  "
Set-TextValue "D14" "98.415.41"
Set-TextValue "E14" "  +0.38%  "
Set-TextValue "D15" "0.0000252"
Set-TextValue "E15" "  +2.07%  "
Set-TextValue "D16" "3.938.46"
Set-TextValue "E16" "  -1.33%  "
Set-TextValue "D17" "5.39"
Set-TextValue "E17" "  -1.89%  "
Set-TextValue "D18" "3.314.08"
Set-TextValue "E18" "  -1.83%  "
Set-TextValue "D19" "15.93"
Set-TextValue "E19" "  +6.97%  "
Set-TextValue "E20" "  -5.04%  "
Set-TextValue "D21" "6.48"
Set-TextValue "E21" "  +9.38%  "
Set-TextValue "D22" "487.60"
Set-TextValue "E22" "  +0.52%  "
Set-TextValue "D23" "9.48"
Set-TextValue "E23" "  +2.25%  "
Set-TextValue "E24" "  -3.53%  "
Set-TextValue "D25" "5.82"
Set-TextValue "E25" "  +1.41%  "
Set-TextValue "D26" "0.347"
Set-TextValue "E26" "  +36.23%  "
Set-TextValue "D27" "89.45"
Set-TextValue "E27" "  +1.41%  "
Set-TextValue "D28" "12.19"
Set-TextValue "E28" "  +1.15%  "
Set-TextValue "D29" "3.491.49"
Set-TextValue "E29" "  -1.48%  "
Set-TextValue "D30" "0.153"
Set-TextValue "E30" "  +22.50%  "
Set-TextValue "E31" "  -0.04%  "
Set-TextValue "E32" "  +2.09%  "
Set-TextValue "D33" "10.71"
Set-TextValue "E33" "  +15.72%  "
Set-TextValue "D34" "1.00"
Set-TextValue "E34" "  +0.10%  "
Set-TextValue "D35" "28.14"
Set-TextValue "E35" "  +2.31%  "
Set-TextValue "D36" "0.484"
Set-TextValue "E36" "  +7.63%  "
Set-TextValue "B37" "RenderToken"
Set-TextValue "C37" "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
Set-TextValue "D37" "7.38"
Set-TextValue "E37" "  +0.11%  "
Set-TextValue "B38" "Kaspa"
Set-TextValue "C38" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D38" "0.151"
Set-TextValue "E38" "  -1.01%  "
Set-TextValue "E39" "  +0.83%  "
Set-TextValue "D40" "498.66"
Set-TextValue "E40" "  -5.18%  "
Set-TextValue "D41" "3.87"
Set-TextValue "E41" "  +2.91%  "
Set-TextValue "E42" "  -0.32%  "
Set-TextValue "E43" "  -1.88%  "
Set-TextValue "E44" "  +0.34%  "
Set-TextValue "E45" "  +0.03%  "
Set-TextValue "D46" "3.16"
Set-TextValue "E46" "  -2.79%  "
Set-TextValue "D47" "161.22"
Set-TextValue "E47" "  +0.22%  "
Set-TextValue "E48" "  +1.40%  "
Set-TextValue "D49" "0.855"
Set-TextValue "E49" "  +7.14%  "
Set-TextValue "D50" "7.35"
Set-TextValue "E50" "  +14.26%  "
Set-TextValue "E51" "  +4.71%  "
